$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) stays formatted as text so that
# numeric-looking strings (e.g. "6.01", "1.00") are not coerced
# into numbers and lose their exact textual representation.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '26.252.92'
$ws.Range("E2").Value = '  +0.64%  '
$ws.Range("D3").Value = '1.605.82'
$ws.Range("E3").Value = '  +0.62%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").Value = '212.75'
$ws.Range("E5").Value = '  +0.35%  '
$ws.Range("E6").Value = '  -0.08%  '
$ws.Range("E7").Value = '  +0.10%  '
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("E9").Value = '  -0.18%  '
$ws.Range("D10").Value = '18.12'
$ws.Range("E10").Value = '  +1.36%  '
$ws.Range("E11").Value = '  -0.48%  '
$ws.Range("D12").Value = '1.830.39'
$ws.Range("E12").Value = '  +0.62%  '
$ws.Range("D13").Value = '1.599.62'
$ws.Range("E13").Value = '  +0.22%  '
$ws.Range("D14").Value = '4.03'
$ws.Range("E14").Value = '  +0.51%  '
$ws.Range("D15").Value = '0.514'
$ws.Range("E15").Value = '  +1.08%  '
$ws.Range("D16").Value = '26.258.72'
$ws.Range("E16").Value = '  +0.67%  '
$ws.Range("D17").Value = '62.15'
$ws.Range("E17").Value = '  +3.02%  '
$ws.Range("E18").Value = '  +1.00%  '
$ws.Range("E19").Value = '  -0.06%  '
$ws.Range("D20").Value = '201.54'
$ws.Range("E20").Value = '  -1.40%  '
$ws.Range("D21").Value = '4.28'
$ws.Range("E21").Value = '  +1.09%  '
$ws.Range("D22").Value = '9.31'
$ws.Range("E22").Value = '  +0.07%  '
$ws.Range("D23").Value = '6.01'
$ws.Range("E23").Value = '  +0.73%  '
$ws.Range("E24").Value = '  +2.50%  '
$ws.Range("D25").Value = '144.78'
$ws.Range("E25").Value = '  +2.05%  '
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  -0.07%  '
$ws.Range("E27").Value = '  -2.86%  '
$ws.Range("D28").Value = '15.22'
$ws.Range("E28").Value = '  +0.08%  '
$ws.Range("D29").Value = '6.56'
$ws.Range("E29").Value = '  +2.04%  '
$ws.Range("D30").Value = '0.0494'
$ws.Range("E30").Value = '  +5.33%  '
$ws.Range("D31").Value = '1.18'
$ws.Range("E31").Value = '  +0.71%  '
$ws.Range("E32").Value = '  +2.77%  '
$ws.Range("E33").Value = '  -1.71%  '
$ws.Range("E34").Value = '  +2.63%  '
$ws.Range("D35").Value = '1.49'
$ws.Range("E35").Value = '  +1.06%  '
$ws.Range("D36").Value = '1.161.14'
$ws.Range("E36").Value = '  +4.97%  '
$ws.Range("D37").Value = '0.0166'
$ws.Range("E37").Value = '  +3.35%  '
$ws.Range("E38").Value = '  -0.05%  '
$ws.Range("D39").Value = '2.33'
$ws.Range("E39").Value = '  +0.32%  '
$ws.Range("D40").Value = '0.789'
$ws.Range("E40").Value = '  +1.67%  '
$ws.Range("E41").Value = '  +0.99%  '
$ws.Range("D42").Value = '0.783'
$ws.Range("E42").Value = '  +0.71%  '
$ws.Range("D44").Value = '1.742.97'
$ws.Range("E44").Value = '  +0.56%  '
$ws.Range("D45").Value = '91.98'
$ws.Range("E45").Value = '  -0.54%  '
$ws.Range("D46").Value = '1.53'
$ws.Range("E46").Value = '  +1.00%  '
$ws.Range("D47").Value = '54.17'
$ws.Range("E47").Value = '  +1.76%  '
$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").Value = '0.0₇0990'
$ws.Range("E48").Value = '  -3.44%  '
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").Value = '0.0505'
$ws.Range("E49").Value = '  +0.08%  '
$ws.Range("E50").Value = '  -0.51%  '
$ws.Range("E51").Value = '  -0.07%  '
